# issue #5: add portion and total(area*portion) of land
#
# The land sheet (sheet1 / "土地") gains two new header columns, "portion"
# and "total" (area * portion), and its single existing data row is
# removed (it only had a share fraction, no computed total yet).
#
# The other four sheets (建物/汽車/現金/存款) have their header row
# collapsed away: the first data row's values are promoted into the
# (bold/bordered) header row slot, and the now-redundant original first
# data row is deleted - subsequent rows shift up keeping their own
# per-row "index" numbers intact.

$wb = $excel.ActiveWorkbook

# ---- Sheet1: 土地 (land) ----------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Delete()
$ws1.Range("P1").Value = "portion"
$ws1.Range("Q1").Value = "total"
$ws1.Range("B1").Copy()
$ws1.Range("P1:Q1").PasteSpecial(-4122)
$ws1.Range("A1").Value = ""

# ---- Sheet2: 建物 (building) -------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2:H2").Copy()
$ws2.Range("B1").PasteSpecial(-4163)
$ws2.Rows.Item(2).Delete()

# ---- Sheet3: 汽車 (car) -------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2:G2").Copy()
$ws3.Range("B1").PasteSpecial(-4163)
$ws3.Rows.Item(2).Delete()

# ---- Sheet4: 現金 (cash) ------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2:D2").Copy()
$ws4.Range("B1").PasteSpecial(-4163)
$ws4.Rows.Item(2).Delete()

# ---- Sheet5: 存款 (deposit) ---------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2:F2").Copy()
$ws5.Range("B1").PasteSpecial(-4163)
$ws5.Rows.Item(2).Delete()

Write-Output "done"
